$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 5.230988666666666
$ws.Range("H2").Value = 15.692966
$ws.Range("I2").Value = 0.2129406655351238
$ws.Range("J2").Value = 0.2129406655351238
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5706193333333334
$ws.Range("N2").Value = 1.711858
$ws.Range("O2").Value = 0.4188640502130462
$ws.Range("P2").Value = 0.4188640502130463
$ws.Range("Q2").Value = 2.984903265647555
$ws.Range("R2").Value = 26.864129390828
$ws.Range("S2").Value = 0.08919318962110356
$ws.Range("T2").Value = 0.08919318962110358
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 5.230988666666666
$ws.Range("H3").Value = 15.692966
$ws.Range("I3").Value = 0.2129406655351238
$ws.Range("J3").Value = 0.2129406655351238
$ws.Range("M3").Value = 0.4846943333333333
$ws.Range("N3").Value = 1.454083
$ws.Range("O3").Value = 0.3557906641356566
$ws.Range("P3").Value = 0.3557906641356566
$ws.Range("Q3").Value = 2.535430564464221
$ws.Range("R3").Value = 22.818875080178
$ws.Range("S3").Value = 0.0757623008122304
$ws.Range("T3").Value = 0.07576230081223041
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 5.230988666666666
$ws.Range("H4").Value = 15.692966
$ws.Range("I4").Value = 0.2129406655351238
$ws.Range("J4").Value = 0.2129406655351238
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3069883333333334
$ws.Range("N4").Value = 0.920965
$ws.Range("O4").Value = 0.2253452856512971
$ws.Range("P4").Value = 0.2253452856512971
$ws.Range("Q4").Value = 1.605852492465556
$ws.Range("R4").Value = 14.45267243219
$ws.Range("S4").Value = 0.04798517510178978
$ws.Range("T4").Value = 0.04798517510178978
$ws.Range("D5").Value = "ECs"
$ws.Range("I5").Value = 0.345577477529236
$ws.Range("J5").Value = 0.3455774775292359
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.5706193333333334
$ws.Range("N5").Value = 1.711858
$ws.Range("O5").Value = 0.4188640502130462
$ws.Range("P5").Value = 0.4188640502130463
$ws.Range("Q5").Value = 4.844144440983333
$ws.Range("R5").Value = 43.59729996885
$ws.Range("S5").Value = 0.1447499819003037
$ws.Range("T5").Value = 0.1447499819003037
$ws.Range("D6").Value = "FAPs"
$ws.Range("I6").Value = 0.345577477529236
$ws.Range("J6").Value = 0.3455774775292359
$ws.Range("M6").Value = 0.4846943333333333
$ws.Range("N6").Value = 1.454083
$ws.Range("O6").Value = 0.3557906641356566
$ws.Range("P6").Value = 0.3557906641356566
$ws.Range("Q6").Value = 4.114703486608333
$ws.Range("R6").Value = 37.03233137947499
$ws.Range("S6").Value = 0.1229532402404518
$ws.Range("T6").Value = 0.1229532402404518
$ws.Range("D7").Value = "MuSCs"
$ws.Range("I7").Value = 0.345577477529236
$ws.Range("J7").Value = 0.3455774775292359
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3069883333333334
$ws.Range("N7").Value = 0.920965
$ws.Range("O7").Value = 0.2253452856512971
$ws.Range("P7").Value = 0.2253452856512971
$ws.Range("Q7").Value = 2.606108383458333
$ws.Range("R7").Value = 23.454975451125
$ws.Range("S7").Value = 0.07787425538848038
$ws.Range("T7").Value = 0.07787425538848038
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 8.418577333333333
$ws.Range("H8").Value = 25.255732
$ws.Range("I8").Value = 0.3426995496362334
$ws.Range("J8").Value = 0.3426995496362334
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.5706193333333334
$ws.Range("N8").Value = 1.711858
$ws.Range("O8").Value = 0.4188640502130462
$ws.Range("P8").Value = 0.4188640502130463
$ws.Range("Q8").Value = 4.803802985561778
$ws.Range("R8").Value = 43.23422687005601
$ws.Range("S8").Value = 0.1435445213668196
$ws.Range("T8").Value = 0.1435445213668196
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 8.418577333333333
$ws.Range("H9").Value = 25.255732
$ws.Range("I9").Value = 0.3426995496362334
$ws.Range("J9").Value = 0.3426995496362334
$ws.Range("M9").Value = 0.4846943333333333
$ws.Range("N9").Value = 1.454083
$ws.Range("O9").Value = 0.3557906641356566
$ws.Range("P9").Value = 0.3557906641356566
$ws.Range("Q9").Value = 4.080436728195111
$ws.Range("R9").Value = 36.723930553756
$ws.Range("S9").Value = 0.1219293003640659
$ws.Range("T9").Value = 0.1219293003640659
$ws.Range("D10").Value = "MuSCs"
$ws.Range("G10").Value = 8.418577333333333
$ws.Range("H10").Value = 25.255732
$ws.Range("I10").Value = 0.3426995496362334
$ws.Range("J10").Value = 0.3426995496362334
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3069883333333334
$ws.Range("N10").Value = 0.920965
$ws.Range("O10").Value = 0.2253452856512971
$ws.Range("P10").Value = 0.2253452856512971
$ws.Range("Q10").Value = 2.584405024597778
$ws.Range("R10").Value = 23.25964522138
$ws.Range("S10").Value = 0.07722572790534789
$ws.Range("T10").Value = 0.07722572790534789
$ws.Range("D11").Value = "ECs"
$ws.Range("G11").Value = 2.426634333333333
$ws.Range("H11").Value = 7.279902999999999
$ws.Range("I11").Value = 0.09878230729940689
$ws.Range("J11").Value = 0.09878230729940687
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.5706193333333334
$ws.Range("N11").Value = 1.711858
$ws.Range("O11").Value = 0.4188640502130462
$ws.Range("P11").Value = 0.4188640502130463
$ws.Range("Q11").Value = 1.384684465530444
$ws.Range("R11").Value = 12.462160189774
$ws.Range("S11").Value = 0.04137635732481933
$ws.Range("T11").Value = 0.04137635732481933
$ws.Range("D12").Value = "FAPs"
$ws.Range("G12").Value = 2.426634333333333
$ws.Range("H12").Value = 7.279902999999999
$ws.Range("I12").Value = 0.09878230729940689
$ws.Range("J12").Value = 0.09878230729940687
$ws.Range("M12").Value = 0.4846943333333333
$ws.Range("N12").Value = 1.454083
$ws.Range("O12").Value = 0.3557906641356566
$ws.Range("P12").Value = 0.3557906641356566
$ws.Range("Q12").Value = 1.176175910438777
$ws.Range("R12").Value = 10.585583193949
$ws.Range("S12").Value = 0.0351458227189085
$ws.Range("T12").Value = 0.0351458227189085
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 2.426634333333333
$ws.Range("H13").Value = 7.279902999999999
$ws.Range("I13").Value = 0.09878230729940689
$ws.Range("J13").Value = 0.09878230729940687
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.3069883333333334
$ws.Range("N13").Value = 0.920965
$ws.Range("O13").Value = 0.2253452856512971
$ws.Range("P13").Value = 0.2253452856512971
$ws.Range("Q13").Value = 0.7449484295994444
$ws.Range("R13").Value = 6.704535866394999
$ws.Range("S13").Value = 0.02226012725567906
$ws.Range("T13").Value = 0.02226012725567905
